# Automatic Excel update [2025-07-28 07:29:18]
# Refresh the "Data ostatniej aktualizacji" (last-update date) column from
# 2025-07-27 to 2025-07-28 for every currently-active offer row on both
# sheets. A couple of duplicate/"hpr" rows on each sheet are intentionally
# left untouched, matching the upstream scraper's dedup behaviour.

$wb = $excel.ActiveWorkbook

$newDate = "2025-07-28"

# --- Sheet 1: "powiat krakowski" ---------------------------------------
$ws1 = $wb.Worksheets.Item("powiat krakowski")
$rng1 = $ws1.Range("E2:E59")
$rng1.NumberFormat = "@"
$rng1.Value = $newDate

# --- Sheet 2: "powiat wielicki" -----------------------------------------
$ws2 = $wb.Worksheets.Item("powiat wielicki")

$ws2.Range("E2:E14").NumberFormat = "@"
$ws2.Range("E2:E14").Value = $newDate

$ws2.Range("E16:E53").NumberFormat = "@"
$ws2.Range("E16:E53").Value = $newDate

$ws2.Range("E55").NumberFormat = "@"
$ws2.Range("E55").Value = $newDate
